# Update the result of week 6
# Adds the "경주 숙소 2박 (8/10-12)" expense line to the 지출 (Expense) sheet,
# recalculates the related summary totals, and adds a new 합계 (Total) block
# on the right-hand side summary table.

$wb = $excel.ActiveWorkbook
$wsIncome  = $wb.Worksheets.Item(1)   # 수입
$wsExpense = $wb.Worksheets.Item(2)   # 지출

# ---- New expense row: 경주 숙소 2박 (8/10-12) on 2025-08-12 (serial 45756) ----
$wsExpense.Range("A11").Value2 = 45756
$wsExpense.Range("B11").Value2 = "경주 숙소 2박 (8/10-12)"
$wsExpense.Range("C11").Value2 = -500480

# ---- Move the running-total row from 11 down to 12 and extend its range ----
$wsExpense.Range("D12").Value2 = "합계"
$wsExpense.Range("E12").Formula = "=SUM(C3:C12)"

# ---- G3 (숙박(한국) 합계) now also includes the new C11 expense ----
$wsExpense.Range("G3").Formula = "=C8+C11"

# ---- Bottom summary formula now points at the relocated total (E12) ----
$wsExpense.Range("E22").Formula = "=SUM(E12,E20)"

# ---- New "수입 유형별 분류"-style mini summary block in F8:G10 ----
$wsExpense.Range("F8").Value2 = "수입"
$wsExpense.Range("G8").Formula = "=수입!E33"

$wsExpense.Range("F9").Value2 = "지출"
$wsExpense.Range("G9").Formula = "=E12"

$wsExpense.Range("F10").Value2 = "합계"
$wsExpense.Range("G10").Formula = "=SUM(G8:G9)"

# ---- Column G width tweak (the old col-level style is no longer needed) ----
$wsExpense.Columns.Item(7).ColumnWidth = 10.875
